# Updates cryptos list prices / hourly volume percentages.
# Note: Price column (D) values that look like plain numbers (e.g. "397.15")
# are prefixed with a leading apostrophe so Excel stores them as literal
# text (matching the workbook's existing inline-string cells) instead of
# silently converting them to floating point numbers. Values that already
# contain multiple "." separators (e.g. "63.544.66") aren't valid numbers
# to begin with, so they don't need the apostrophe.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.544.66"
$ws.Range("E2").Value = "  -6.49%  "

$ws.Range("D3").Value = "3.561.54"
$ws.Range("E3").Value = "  -4.40%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'397.15"
$ws.Range("E5").Value = "  -6.62%  "

$ws.Range("D6").Value = "'122.91"
$ws.Range("E6").Value = "  -6.27%  "

$ws.Range("D7").Value = "3.558.75"
$ws.Range("E7").Value = "  -4.27%  "

$ws.Range("E8").Value = "  -9.56%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "'0.681"
$ws.Range("E10").Value = "  -11.95%  "

$ws.Range("E11").Value = "  -18.98%  "

$ws.Range("D12").Value = "'0.0000326"
$ws.Range("E12").Value = "  -20.46%  "

$ws.Range("D13").Value = "'38.95"
$ws.Range("E13").Value = "  -9.26%  "

$ws.Range("D14").Value = "4.121.23"
$ws.Range("E14").Value = "  -4.22%  "

$ws.Range("D15").Value = "'9.19"
$ws.Range("E15").Value = "  -7.53%  "

$ws.Range("E16").Value = "  -2.83%  "

$ws.Range("D17").Value = "'13.88"
$ws.Range("E17").Value = "  +7.81%  "

$ws.Range("D18").Value = "3.568.35"
$ws.Range("E18").Value = "  -3.84%  "

$ws.Range("D19").Value = "'18.71"
$ws.Range("E19").Value = "  -9.20%  "

$ws.Range("D20").Value = "63.593.37"
$ws.Range("E20").Value = "  -6.44%  "

$ws.Range("E21").Value = "  -10.28%  "

$ws.Range("D22").Value = "'394.51"
$ws.Range("E22").Value = "  -12.79%  "

$ws.Range("D23").Value = "'13.92"
$ws.Range("E23").Value = "  -8.17%  "

$ws.Range("D24").Value = "'81.99"
$ws.Range("E24").Value = "  -8.77%  "

$ws.Range("D25").Value = "'2.93"
$ws.Range("E25").Value = "  -6.62%  "

$ws.Range("E26").Value = "  +8.83%  "

$ws.Range("D27").Value = "'33.97"
$ws.Range("E27").Value = "  -11.97%  "

$ws.Range("E28").Value = "  -8.78%  "

$ws.Range("D29").Value = "'8.74"
$ws.Range("E29").Value = "  -16.07%  "

$ws.Range("D30").Value = "'11.96"
$ws.Range("E30").Value = "  -4.14%  "

$ws.Range("D31").Value = "'2.60"
$ws.Range("E31").Value = "  -7.80%  "

$ws.Range("E32").Value = "  -7.05%  "

$ws.Range("D33").Value = "'6.83"
$ws.Range("E33").Value = "  -5.54%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("E36").Value = "  -9.75%  "

$ws.Range("D37").Value = "'54.01"
$ws.Range("E37").Value = "  -4.50%  "

$ws.Range("E38").Value = "  -11.26%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "0.0₃0652"
$ws.Range("E40").Value = "  -11.92%  "

$ws.Range("E41").Value = "  -10.95%  "

$ws.Range("D42").Value = "'0.131"
$ws.Range("E42").Value = "  -11.40%  "

$ws.Range("E43").Value = "  +16.42%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'26.09"
$ws.Range("E44").Value = "  -2.43%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'142.05"
$ws.Range("E45").Value = "  -3.34%  "

$ws.Range("E46").Value = "  -6.22%  "

$ws.Range("E47").Value = "  -10.43%  "

$ws.Range("E48").Value = "  -7.60%  "

$ws.Range("E49").Value = "  -6.50%  "

$ws.Range("E50").Value = "  -9.17%  "

$ws.Range("D51").Value = "'0.277"
$ws.Range("E51").Value = "  -9.82%  "
